$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old content
$ws.Range("A1:C6").Clear()

# Headers
$ws.Range("B1").Value = "Task"
$ws.Range("E1").Value = "Assignee"
$ws.Range("G1").Value = "Due Date"
$ws.Range("G1").Font.Color = 0
$ws.Range("I1").Value = "Status"

# Assignee column
$ws.Range("E2").Value = "Ganden"
$ws.Range("E3").Value = "Marco"

# Task column
$ws.Range("B2").Value = "Laundry"
$ws.Range("B3").Value = "Forge signatures"
$ws.Range("B4").Value = "Snapper"
$ws.Range("B5").Value = "Gapper"
$ws.Range("B6").Value = "Siege"
$ws.Range("B7").Value = "Impact"

# Due date column
$ws.Range("G2").NumberFormat = "mm-dd-yy"
$ws.Range("G2").Value = 45266
$ws.Range("G2").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("G3").Value = 45278

# Status column
$ws.Range("I2").Value = $false
$ws.Range("I3").Value = $true

# Selection
$ws.Range("B7").Select()
